$wb = $excel.ActiveWorkbook

# --- Sheet 1: Summary ---
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.9537366548042705
$ws1.Range("C2").Value = 0.5416666666666666
$ws1.Range("D2").Value = 0.4642857142857143
$ws1.Range("E2").Value = 0.5
$ws1.Range("F2").Value = 0.4779411764705883
$ws1.Range("G2").Value = 0.4668508287292817
$ws1.Range("H2").Value = 0.7218432316746924
$ws1.Range("I2").Value = 13
$ws1.Range("J2").Value = 11
$ws1.Range("K2").Value = 523
$ws1.Range("L2").Value = 15

# --- Sheet 2: Classification Report ---
$ws2 = $wb.Worksheets.Item("Classification Report")
$ws2.Range("B2").Value = 0.9721189591078067
$ws2.Range("C2").Value = 0.9794007490636704
$ws2.Range("D2").Value = 0.9757462686567164

$ws2.Range("B3").Value = 0.5416666666666666
$ws2.Range("C3").Value = 0.4642857142857143
$ws2.Range("D3").Value = 0.5

$ws2.Range("B4").Value = 0.9537366548042705
$ws2.Range("C4").Value = 0.9537366548042705
$ws2.Range("D4").Value = 0.9537366548042705
$ws2.Range("E4").Value = 0.9537366548042705

$ws2.Range("B5").Value = 0.7568928128872366
$ws2.Range("C5").Value = 0.7218432316746923
$ws2.Range("D5").Value = 0.7378731343283582

$ws2.Range("B6").Value = 0.9506729374203478
$ws2.Range("C6").Value = 0.9537366548042705
$ws2.Range("D6").Value = 0.9520436075848515

# --- Sheet 3: Confusion Matrix ---
$ws3 = $wb.Worksheets.Item("Confusion Matrix")
$ws3.Range("B2").Value = 523
$ws3.Range("C2").Value = 11
$ws3.Range("B3").Value = 15
$ws3.Range("C3").Value = 13
